# Refresh KHL injuries snapshot: mark "Гараев Амир" (СОЧ) as returned from injury,
# add newly-injured "Уотерспун Тайлер" (БАР) into the active snapshot, and
# refresh the scraped_at timestamps for every row re-checked in this run.

$wb = $excel.ActiveWorkbook

$wsSnapshot   = $wb.Worksheets.Item("snapshot")
$wsReturned   = $wb.Worksheets.Item("returned")
$wsNewInjured = $wb.Worksheets.Item("new_injured")

# Full replacement data for the "snapshot" sheet, rows 2-48 (row 1 is the header
# and is left untouched). This already reflects: the new БАР/Уотерспун Тайлер
# injury inserted in alphabetical order, the removal of СОЧ/Гараев Амир (who
# returned from injury), and refreshed scraped_at values for every row.
$dataSnapshot = @(
  @("АВГ","Авангард","avangard","Якупов Наиль","65","нападающий","16391","1369_АВГ_якуповнаиль","injured_active","https://www.khl.ru/clubs/avangard/team/","2025-11-04T09:22:56.615985+00:00"),
  @("АВТ","Автомобилист","avtomobilist","Зборовский Сергей","2","защитник","20989","1369_АВТ_зборовскийсергей","injured_active","https://www.khl.ru/clubs/avtomobilist/team/","2025-11-04T09:22:59.139224+00:00"),
  @("АВТ","Автомобилист","avtomobilist","Кизимов Семён","97","нападающий","25697","1369_АВТ_кизимовсемен","injured_active","https://www.khl.ru/clubs/avtomobilist/team/","2025-11-04T09:22:59.139245+00:00"),
  @("АВТ","Автомобилист","avtomobilist","Осипов Максим И","21","защитник","17459","1369_АВТ_осиповмаксими","injured_active","https://www.khl.ru/clubs/avtomobilist/team/","2025-11-04T09:22:59.139272+00:00"),
  @("АВТ","Автомобилист","avtomobilist","Трямкин Никита","88","защитник","17594","1369_АВТ_трямкинникита","injured_active","https://www.khl.ru/clubs/avtomobilist/team/","2025-11-04T09:22:59.139280+00:00"),
  @("АДМ","Адмирал","admiral","Грман Марио","77","защитник","31232","1369_АДМ_грманмарио","injured_active","https://www.khl.ru/clubs/admiral/team/","2025-11-04T09:23:01.216155+00:00"),
  @("АДМ","Адмирал","admiral","Старков Степан","18","нападающий","27000","1369_АДМ_старковстепан","injured_active","https://www.khl.ru/clubs/admiral/team/","2025-11-04T09:23:01.216195+00:00"),
  @("АДМ","Адмирал","admiral","Шепелев Александр","19","защитник","23447","1369_АДМ_шепелевалександр","injured_active","https://www.khl.ru/clubs/admiral/team/","2025-11-04T09:23:01.216219+00:00"),
  @("АКБ","Ак Барс","ak_bars","Яруллин Альберт","33","защитник","16365","1369_АКБ_яруллинальберт","injured_active","https://www.khl.ru/clubs/ak_bars/team/","2025-11-04T09:23:03.917959+00:00"),
  @("АМР","Амур","amur","Абросимов Роман","94","защитник","17968","1369_АМР_абросимовроман","injured_active","https://www.khl.ru/clubs/amur/team/","2025-11-04T09:23:06.320029+00:00"),
  @("АМР","Амур","amur","Броадхёрст Алекс","25","нападающий","27232","1369_АМР_броадхерсталекс","injured_active","https://www.khl.ru/clubs/amur/team/","2025-11-04T09:23:06.320060+00:00"),
  @("АМР","Амур","amur","Гиздатуллин Артур","87","нападающий","22208","1369_АМР_гиздатуллинартур","injured_active","https://www.khl.ru/clubs/amur/team/","2025-11-04T09:23:06.320079+00:00"),
  @("БАР","Барыс","barys","Бояркин Никита","1","вратарь","28244","1369_БАР_бояркинникита","injured_active","https://www.khl.ru/clubs/barys/team/","2025-11-04T09:23:08.784893+00:00"),
  @("БАР","Барыс","barys","Галимов Эмиль","27","нападающий","15997","1369_БАР_галимовэмиль","injured_active","https://www.khl.ru/clubs/barys/team/","2025-11-04T09:23:08.784924+00:00"),
  @("БАР","Барыс","barys","Мухаметов Максим","23","нападающий","25207","1369_БАР_мухаметовмаксим","injured_active","https://www.khl.ru/clubs/barys/team/","2025-11-04T09:23:08.784944+00:00"),
  @("БАР","Барыс","barys","Уотерспун Тайлер","26","защитник","45769","1369_БАР_уотерспунтайлер","injured_active","https://www.khl.ru/clubs/barys/team/","2025-11-04T09:23:08.784962+00:00"),
  @("ЛАД","Лада","lada","Ожгихин Алексей","43","нападающий","23021","1369_ЛАД_ожгихиналексей","injured_active","https://www.khl.ru/clubs/lada/team/","2025-11-04T09:23:16.357275+00:00"),
  @("ЛОК","Локомотив","lokomotiv","Сергеев Андрей","99","защитник","15416","1369_ЛОК_сергеевандрей","injured_active","https://www.khl.ru/clubs/lokomotiv/team/","2025-11-04T09:23:18.733981+00:00"),
  @("ММГ","Металлург Мг","metallurg_mg","Козлов Андрей Е","39","нападающий","40899","1369_ММГ_козловандрейе","injured_active","https://www.khl.ru/clubs/metallurg_mg/team/","2025-11-04T09:23:21.645334+00:00"),
  @("НХК","Нефтехимик","neftekhimik","Дергачёв Александр","18","нападающий","20592","1369_НХК_дергачевалександр","injured_active","https://www.khl.ru/clubs/neftekhimik/team/","2025-11-04T09:23:24.489597+00:00"),
  @("НХК","Нефтехимик","neftekhimik","Попугаев Никита О","13","нападающий","22683","1369_НХК_попугаевникитао","injured_active","https://www.khl.ru/clubs/neftekhimik/team/","2025-11-04T09:23:24.489629+00:00"),
  @("НХК","Нефтехимик","neftekhimik","Профака Лука","22","защитник","43943","1369_НХК_профакалука","injured_active","https://www.khl.ru/clubs/neftekhimik/team/","2025-11-04T09:23:24.489647+00:00"),
  @("СЕВ","Северсталь","severstal","Ващенко Григорий","16","защитник","14155","1369_СЕВ_ващенкогригорий","injured_active","https://www.khl.ru/clubs/severstal/team/","2025-11-04T09:23:27.113355+00:00"),
  @("СЕВ","Северсталь","severstal","Грудинин Владимир","2","защитник","35064","1369_СЕВ_грудининвладимир","injured_active","https://www.khl.ru/clubs/severstal/team/","2025-11-04T09:23:27.113386+00:00"),
  @("СЕВ","Северсталь","severstal","Танков Кирилл","42","нападающий","32981","1369_СЕВ_танковкирилл","injured_active","https://www.khl.ru/clubs/severstal/team/","2025-11-04T09:23:27.113406+00:00"),
  @("СЕВ","Северсталь","severstal","Цицюра Владислав","10","нападающий","23840","1369_СЕВ_цицюравладислав","injured_active","https://www.khl.ru/clubs/severstal/team/","2025-11-04T09:23:27.113424+00:00"),
  @("СИБ","Сибирь","sibir","Приски Чейз Эванс","22","защитник","45392","1369_СИБ_прискичейзэванс","injured_active","https://www.khl.ru/clubs/sibir/team/","2025-11-04T09:23:29.222616+00:00"),
  @("СОЧ","ХК Сочи","hc_sochi","Гуськов Матвей","77","нападающий","29136","1369_СОЧ_гуськовматвей","injured_active","https://www.khl.ru/clubs/hc_sochi/team/","2025-11-04T09:23:33.391115+00:00"),
  @("СОЧ","ХК Сочи","hc_sochi","Мачулин Василий","23","защитник","33926","1369_СОЧ_мачулинвасилий","injured_active","https://www.khl.ru/clubs/hc_sochi/team/","2025-11-04T09:23:33.391147+00:00"),
  @("СОЧ","ХК Сочи","hc_sochi","Хомченко Павел","30","вратарь","17592","1369_СОЧ_хомченкопавел","injured_active","https://www.khl.ru/clubs/hc_sochi/team/","2025-11-04T09:23:33.391167+00:00"),
  @("СПР","Спартак","spartak","Воробьёв Иван В","88","нападающий","33545","1369_СПР_воробьевиванв","injured_active","https://www.khl.ru/clubs/spartak/team/","2025-11-04T09:23:35.898391+00:00"),
  @("СПР","Спартак","spartak","Порядин Павел","24","нападающий","19258","1369_СПР_порядинпавел","injured_active","https://www.khl.ru/clubs/spartak/team/","2025-11-04T09:23:35.898421+00:00"),
  @("СПР","Спартак","spartak","Рубцов Герман","95","нападающий","22494","1369_СПР_рубцовгерман","injured_active","https://www.khl.ru/clubs/spartak/team/","2025-11-04T09:23:35.898440+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Алалыкин Данил","61","нападающий","34493","1369_СЮЛ_алалыкинданил","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075626+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Берлёв Антон","83","нападающий","20546","1369_СЮЛ_берлевантон","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075661+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Зоркин Никита","52","защитник","26738","1369_СЮЛ_зоркинникита","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075675+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Кузьмин Глеб","17","нападающий","22170","1369_СЮЛ_кузьминглеб","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075688+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Пименов Артём","68","нападающий","21205","1369_СЮЛ_пименовартем","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075701+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Хворов Николай","44","нападающий","39828","1369_СЮЛ_хворовниколай","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075714+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Хохряков Пётр","62","нападающий","15413","1369_СЮЛ_хохряковпетр","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075734+00:00"),
  @("СЮЛ","Салават Юлаев","salavat_yulaev","Ян Денис","77","нападающий","22166","1369_СЮЛ_янденис","injured_active","https://www.khl.ru/clubs/salavat_yulaev/team/","2025-11-04T09:23:38.075745+00:00"),
  @("ТОР","Торпедо","torpedo","Науменков Михаил","33","защитник","16400","1369_ТОР_науменковмихаил","injured_active","https://www.khl.ru/clubs/torpedo/team/","2025-11-04T09:23:40.181569+00:00"),
  @("ТОР","Торпедо","torpedo","Рожков Никита А","71","нападающий","27912","1369_ТОР_рожковникитаа","injured_active","https://www.khl.ru/clubs/torpedo/team/","2025-11-04T09:23:40.181602+00:00"),
  @("ЦСК","ЦСКА","cska","Моисеев Данила","93","нападающий","23931","1369_ЦСК_моисеевданила","injured_active","https://www.khl.ru/clubs/cska/team/","2025-11-04T09:23:44.806068+00:00"),
  @("ШДР","Драконы","kunlun","Гроло Жереми","75","защитник","45343","1369_ШДР_гроложереми","injured_active","https://www.khl.ru/clubs/kunlun/team/","2025-11-04T09:23:47.510765+00:00"),
  @("ШДР","Драконы","kunlun","Саттер Райли","14","нападающий","45491","1369_ШДР_саттеррайли","injured_active","https://www.khl.ru/clubs/kunlun/team/","2025-11-04T09:23:47.510794+00:00"),
  @("ШДР","Драконы","kunlun","Фу Спенсер","15","нападающий","34934","1369_ШДР_фуспенсер","injured_active","https://www.khl.ru/clubs/kunlun/team/","2025-11-04T09:23:47.510812+00:00")
)

for ($i = 0; $i -lt $dataSnapshot.Length; $i++) {
    $rowValues = $dataSnapshot[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $wsSnapshot.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}

# Append the "returned" log entry for Гараев Амир (СОЧ).
$dataReturned = @(
  @("СОЧ","ХК Сочи","Гараев Амир","1369_СОЧ_гараевамир","RETURN","2025-11-04T17:23:48.013752+08:00","2025-11-04")
)

for ($i = 0; $i -lt $dataReturned.Length; $i++) {
    $rowValues = $dataReturned[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $wsReturned.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}

# Append the "new_injured" log entry for Уотерспун Тайлер (БАР).
$dataNewInjured = @(
  @("БАР","Барыс","Уотерспун Тайлер","1369_БАР_уотерспунтайлер","INJURED_NEW","2025-11-04T17:23:48.013752+08:00","2025-11-04")
)

for ($i = 0; $i -lt $dataNewInjured.Length; $i++) {
    $rowValues = $dataNewInjured[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $wsNewInjured.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
